$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6528
$ws.Range("E2").Value = 792
$ws.Range("F2").Value = 792
$ws.Range("G2").Value = 690
$ws.Range("H2").Value = 552
$ws.Range("I2").Value = 550
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 6830
$ws.Range("L2").Value = 2463
$ws.Range("M2").Value = 4367
$ws.Range("N2").Value = 4328
$ws.Range("O2").Value = 39
$ws.Range("P2").Value = 1196
$ws.Range("Q2").Value = 1401
$ws.Range("R2").Value = -209
$ws.Range("S2").Value = -264
$ws.Range("T2").Value = 946
$ws.Range("U2").Value = 455
$ws.Range("V2").Value = 837
$ws.Range("W2").Value = 12.13
$ws.Range("X2").Value = 8.45
$ws.Range("Y2").Value = 13.14
$ws.Range("Z2").Value = 8.08
$ws.Range("AA2").Value = 56.41
$ws.Range("AB2").Value = 267.64
$ws.Range("AC2").Value = 1151
$ws.Range("AD2").Value = 15.94
$ws.Range("AE2").Value = 9100
$ws.Range("AF2").Value = 2.02
$ws.Range("AG2").Value = 350
$ws.Range("AH2").Value = 1.91
$ws.Range("AI2").Value = 30.24
$ws.Range("AJ2").Value = 47821966

# Row 3
$ws.Range("D3").Value = 6610
$ws.Range("E3").Value = 996
$ws.Range("F3").Value = 996
$ws.Range("G3").Value = 971
$ws.Range("H3").Value = 730
$ws.Range("I3").Value = 727
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 7113
$ws.Range("L3").Value = 2178
$ws.Range("M3").Value = 4934
$ws.Range("N3").Value = 4893
$ws.Range("O3").Value = 42
$ws.Range("P3").Value = 1196
$ws.Range("Q3").Value = 1578
$ws.Range("R3").Value = -924
$ws.Range("S3").Value = -360
$ws.Range("T3").Value = 845
$ws.Range("U3").Value = 733
$ws.Range("V3").Value = 646
$ws.Range("W3").Value = 15.07
$ws.Range("X3").Value = 11.04
$ws.Range("Y3").Value = 15.77
$ws.Range("Z3").Value = 10.47
$ws.Range("AA3").Value = 44.15
$ws.Range("AB3").Value = 314.75
$ws.Range("AC3").Value = 1521
$ws.Range("AD3").Value = 11.38
$ws.Range("AE3").Value = 10287
$ws.Range("AF3").Value = 1.68
$ws.Range("AG3").Value = 350
$ws.Range("AH3").Value = 2.02
$ws.Range("AI3").Value = 22.89
$ws.Range("AJ3").Value = 47821966

# Row 4
$ws.Range("D4").Value = 6651
$ws.Range("E4").Value = 806
$ws.Range("F4").Value = 806
$ws.Range("G4").Value = 813
$ws.Range("H4").Value = 689
$ws.Range("I4").Value = 687
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 7779
$ws.Range("L4").Value = 2315
$ws.Range("M4").Value = 5465
$ws.Range("N4").Value = 5419
$ws.Range("O4").Value = 46
$ws.Range("P4").Value = 1196
$ws.Range("Q4").Value = 1554
$ws.Range("R4").Value = -2105
$ws.Range("S4").Value = -166
$ws.Range("T4").Value = 712
$ws.Range("U4").Value = 842
$ws.Range("V4").Value = 639
$ws.Range("W4").Value = 12.13
$ws.Range("X4").Value = 10.36
$ws.Range("Y4").Value = 13.32
$ws.Range("Z4").Value = 9.25
$ws.Range("AA4").Value = 42.35
$ws.Range("AB4").Value = 358.83
$ws.Range("AC4").Value = 1436
$ws.Range("AD4").Value = 12.05
$ws.Range("AE4").Value = 11393
$ws.Range("AF4").Value = 1.52
$ws.Range("AG4").Value = 415
$ws.Range("AH4").Value = 2.4
$ws.Range("AI4").Value = 28.75
$ws.Range("AJ4").Value = 47821966

# Row 5
$ws.Range("D5").Value = 6858
$ws.Range("E5").Value = 754
$ws.Range("F5").Value = 754
$ws.Range("G5").Value = 749
$ws.Range("H5").Value = 573
$ws.Range("I5").Value = 574
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 7929
$ws.Range("L5").Value = 2106
$ws.Range("M5").Value = 5823
$ws.Range("N5").Value = 5778
$ws.Range("O5").Value = 45
$ws.Range("P5").Value = 1196
$ws.Range("Q5").Value = 993
$ws.Range("R5").Value = -818
$ws.Range("S5").Value = -197
$ws.Range("T5").Value = 579
$ws.Range("U5").Value = 414
$ws.Range("V5").Value = 641
$ws.Range("W5").Value = 10.99
$ws.Range("X5").Value = 8.36
$ws.Range("Y5").Value = 10.25
$ws.Range("Z5").Value = 7.3
$ws.Range("AA5").Value = 36.16
$ws.Range("AB5").Value = 389.1
$ws.Range("AC5").Value = 1200
$ws.Range("AD5").Value = 11.21
$ws.Range("AE5").Value = 12148
$ws.Range("AF5").Value = 1.11
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 2.6
$ws.Range("AI5").Value = 29.01
$ws.Range("AJ5").Value = 47821966

# Row 6
$ws.Range("D6").Value = 6908
$ws.Range("E6").Value = 667
$ws.Range("F6").Value = 667
$ws.Range("G6").Value = 680
$ws.Range("H6").Value = 520
$ws.Range("I6").Value = 516
$ws.Range("K6").Value = 8160
$ws.Range("L6").Value = 1498
$ws.Range("M6").Value = 6662
$ws.Range("N6").Value = 6613
$ws.Range("P6").Value = 1196
$ws.Range("Q6").Value = 1835
$ws.Range("R6").Value = -1398
$ws.Range("S6").Value = -776
$ws.Range("T6").Value = 686
$ws.Range("U6").Value = 1149
$ws.Range("V6").Value = 33
$ws.Range("W6").Value = 9.66
$ws.Range("X6").Value = 7.53
$ws.Range("Y6").Value = 8.33
$ws.Range("Z6").Value = 6.46
$ws.Range("AA6").Value = 22.49
$ws.Range("AB6").Value = 458.83
$ws.Range("AC6").Value = 1079
$ws.Range("AD6").Value = 10.65
$ws.Range("AE6").Value = 13903
$ws.Range("AF6").Value = 0.83
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 3.04
$ws.Range("AI6").Value = 32.25
$ws.Range("AJ6").Value = 47821966

# Row 7
$ws.Range("D7").Value = 7114
$ws.Range("E7").Value = 672
$ws.Range("G7").Value = 737
$ws.Range("H7").Value = 540
$ws.Range("I7").Value = 532
$ws.Range("K7").Value = 8564
$ws.Range("L7").Value = 1567
$ws.Range("M7").Value = 6997
$ws.Range("N7").Value = 6943
$ws.Range("P7").Value = 1199
$ws.Range("Q7").Value = 1324
$ws.Range("R7").Value = -915
$ws.Range("S7").Value = -192
$ws.Range("T7").Value = 538
$ws.Range("U7").Value = 780
$ws.Range("W7").Value = 9.45
$ws.Range("X7").Value = 7.59
$ws.Range("Y7").Value = 7.84
$ws.Range("Z7").Value = 6.45
$ws.Range("AA7").Value = 22.4
$ws.Range("AC7").Value = 1112
$ws.Range("AD7").Value = 7.65
$ws.Range("AE7").Value = 14597
$ws.Range("AF7").Value = 0.58
$ws.Range("AG7").Value = 366
$ws.Range("AH7").Value = 4.3
$ws.Range("AI7").Value = 32.88

# Row 8
$ws.Range("D8").Value = 7176
$ws.Range("E8").Value = 686
$ws.Range("G8").Value = 694
$ws.Range("H8").Value = 530
$ws.Range("I8").Value = 525
$ws.Range("K8").Value = 8767
$ws.Range("L8").Value = 1430
$ws.Range("M8").Value = 7334
$ws.Range("N8").Value = 7271
$ws.Range("P8").Value = 1199
$ws.Range("Q8").Value = 1133
$ws.Range("R8").Value = -870
$ws.Range("S8").Value = -189
$ws.Range("T8").Value = 525
$ws.Range("U8").Value = 775
$ws.Range("W8").Value = 9.56
$ws.Range("X8").Value = 7.38
$ws.Range("Y8").Value = 7.38
$ws.Range("Z8").Value = 6.11
$ws.Range("AA8").Value = 19.49
$ws.Range("AC8").Value = 1097
$ws.Range("AD8").Value = 7.9
$ws.Range("AE8").Value = 15288
$ws.Range("AF8").Value = 0.57
$ws.Range("AG8").Value = 377
$ws.Range("AH8").Value = 4.34
$ws.Range("AI8").Value = 34.33

# Row 9
$ws.Range("D9").Value = 7265
$ws.Range("E9").Value = 650
$ws.Range("G9").Value = 680
$ws.Range("H9").Value = 520
$ws.Range("I9").Value = 515
$ws.Range("K9").Value = 9085
$ws.Range("L9").Value = 1410
$ws.Range("M9").Value = 7675
$ws.Range("N9").Value = 7605
$ws.Range("P9").Value = 1200
$ws.Range("Q9").Value = 1175
$ws.Range("R9").Value = -910
$ws.Range("S9").Value = -210
$ws.Range("T9").Value = 630
$ws.Range("W9").Value = 8.95
$ws.Range("X9").Value = 7.16
$ws.Range("Y9").Value = 6.92
$ws.Range("Z9").Value = 5.83
$ws.Range("AA9").Value = 18.37
$ws.Range("AC9").Value = 1077
$ws.Range("AD9").Value = 8.05
$ws.Range("AE9").Value = 15989
$ws.Range("AF9").Value = 0.54
$ws.Range("AG9").Value = 395
$ws.Range("AH9").Value = 4.56
$ws.Range("AI9").Value = 36.68

# U9 cell is removed entirely in the target (value deleted)
$ws.Range("U9").ClearContents()
